$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record columns, matching the
# existing header formatting (style used by A1:AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-54) with the
# team's season record: 95 wins, 67 losses, 0 ties.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 95
    $ws.Cells.Item($row, 31).Value = 67
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-54"
